# Paper revision: ontology changes + adding NCBI sample identifiers.
# A new "Open sea and tidal areas" sub-row is inserted right after row 226
# (mfd_hab1 "Open sea and tidal areas" now gets its own standalone row,
# ahead of its mfd_hab2 children), pushing every following ontology row
# down by one and appending one extra trailing row to the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the last (always-blank) worksheet row first so that shifting rows
# down below doesn't push a row past Excel's maximum row (1,048,576).
$ws.Rows("1048576:1048576").Delete()

# Insert a new blank row at 227; rows 227-262 (and the trailing blank
# spacer rows) shift down to 228-263.
$ws.Rows("227:227").Insert()

# Row 226 keeps its mfd_sampletype/areatype/hab1_code/hab1 (A-D) and
# mfd_hab3 group (K) values, and row 227 inherits the same A-D/K values...
$ws.Range("A226:D226").Copy($ws.Range("A227:D227"))
$ws.Range("K226").Copy($ws.Range("K227"))

# ...while the mfd_hab2_code/mfd_hab2 values (E-F) that used to live on
# row 226 move down onto the new row 227, leaving row 226's E-F blank.
$ws.Range("E226:F226").Cut($ws.Range("E227:F227"))

# Restore the active-cell selection recorded in the saved workbook.
$ws.Range("E226").Select()
